# Update the "Test Cases" sheet: change column D (rows 2-22) from "N" to "Y"
# and update the sheet's selection to D2:D24 with D2 active.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

for ($row = 2; $row -le 22; $row++) {
    $ws.Cells.Item($row, 4).Value = "Y"
}

$ws.Activate()
$ws.Range("D2:D24").Select()
$excel.ActiveWindow.ActiveCell = $ws.Range("D2")
